$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text must be pre-formatted as Text so that
# Excel keeps them as strings (matching the original report formatting) instead
# of auto-converting them to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.163.20"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.955.77"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "587.17"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "148.06"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "2.938.80"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "6.80"
$ws.Range("E10").Value = "  +10.80%  "
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "34.51"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "3.446.63"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "6.88"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "61.141.43"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "2.949.86"
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("D20").Value = "433.38"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "0.676"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").Value = "7.31"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "80.48"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "11.02"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "11.91"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.64"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "26.96"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "0.0₃0840"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "5.72"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "50.11"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").Value = "8.77"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  +8.18%  "
$ws.Range("D44").Value = "42.60"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("D45").Value = "0.0349"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").Value = "374.81"
$ws.Range("E46").Value = "  -6.12%  "
$ws.Range("D47").Value = "2.667.16"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "133.52"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "25.84"
$ws.Range("E49").Value = "  +9.36%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "2.16"
$ws.Range("E51").Value = "  -1.09%  "
